$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the old "Locked"/"Lockless"
# label column and all the iteration-count columns one to the right.
$ws.Columns("B:B").Insert()

# New header row: Architecture (B1) and Iterations (D1, above the numeric run).
$ws.Range("B1").Value = "Architecture"
$ws.Range("D1").Value = "Iterations"

# Row 2 iteration counts are now plain numbers instead of "N Iterations" text.
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 10000
$ws.Range("G2").Value = 100000
$ws.Range("H2").Value = 1000000

# New architecture column value for the "Locked" results row.
$ws.Range("B4").Value = "64 Core Intel Xeon CPU E7-4820 @ 2.00GHz (Stoker)"

# Widen column A (labels) and column B (new architecture column) to fit text.
$ws.Columns("A:A").ColumnWidth = 31.666666666666668
$ws.Columns("B:B").ColumnWidth = 45.666666666666664

# Update the selection to reflect the new layout.
$ws.Range("B21").Select()
